# Update "想去人数" (F column) values on both "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    7  = 84
    8  = 454
    11 = 566
    13 = 298
    15 = 369
    19 = 50
    22 = 917
    23 = 1394
    24 = 296
    32 = 246
    33 = 272
    37 = 156
    38 = 579
    40 = 3617
    42 = 195
    44 = 41
    46 = 64
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
